$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "Resultado real" (column F) for rows 3-6 with the actual results.
$ws.Range("F3").Value = $ws.Range("E3").Value2
$ws.Range("F4").Value = $ws.Range("E4").Value2
$ws.Range("F5").Value = $ws.Range("E5").Value2
$ws.Range("F6").Value = $ws.Range("E5").Value2

# Fill in "Estado" (column G) with "Pasó" for rows 2-6.
$ws.Range("G2").Value = "Pasó"
$ws.Range("G3").Value = "Pasó"
$ws.Range("G4").Value = "Pasó"
$ws.Range("G5").Value = "Pasó"
$ws.Range("G6").Value = "Pasó"

# Remove row 7 (TC06), which was never executed.
$ws.Rows("7").Delete()

# Update the selected cell to match the last edited cell.
$ws.Range("F6").Select()
